$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.11928112952215741
$ws.Range("B1").Value = 0.119242864003553
$ws.Range("A2").Value = -0.097141087933028558
$ws.Range("B2").Value = 0.096992786729967584
$ws.Range("A3").Value = -0.047293486134936558
$ws.Range("B3").Value = 0.047119960398331884
$ws.Range("A4").Value = -0.039119960488728012
$ws.Range("B4").Value = 0.038750951723471871
$ws.Range("A5").Value = -0.035750951767269612
$ws.Range("B5").Value = 0.034500072802053694
$ws.Range("A6").Value = -0.0053143369475368019
$ws.Range("B6").Value = 0.005196709276720668
$ws.Range("A7").Value = 0.0048032906040362811
$ws.Range("B7").Value = -0.0048207829033337646
$ws.Range("A8").Value = 0.014820782784589426
$ws.Range("B8").Value = -0.014840497853476187
$ws.Range("A9").Value = 0.016840497815615585
$ws.Range("B9").Value = -0.01685686629088945
$ws.Range("A10").Value = 0.018856866255715587
$ws.Range("B10").Value = -0.018856569206588958
$ws.Range("A11").Value = 0.02185656916169787
$ws.Range("B11").Value = -0.021859489320403824
$ws.Range("A12").Value = 0.025359489271229663
$ws.Range("B12").Value = -0.025412571290631458
$ws.Range("A13").Value = 0.028912571247679431
$ws.Range("B13").Value = -0.028959634082012187
$ws.Range("A14").Value = 0.0047987551656065719
$ws.Range("B14").Value = -0.0048334563124319274
$ws.Range("A15").Value = 0.0058334562969131198
$ws.Range("B15").Value = -0.0058578693419120142
$ws.Range("A16").Value = 0.007857869317280386
$ws.Range("B16").Value = -0.0079097704303743299
$ws.Range("A17").Value = 0.0099097704076127613
$ws.Range("B17").Value = -0.009920030712113892
$ws.Range("A18").Value = -0.016101776130398093
$ws.Range("B18").Value = 0.016090526894188883
$ws.Range("A19").Value = -0.012090526934544599
$ws.Range("B19").Value = 0.012015855256806152
$ws.Range("A20").Value = -0.0080158553003233379
$ws.Range("B20").Value = 0.0080055340170144262
$ws.Range("A21").Value = -0.0040055340610187784
$ws.Range("B21").Value = 0.0039999999556137311
$ws.Range("A22").Value = -0.045699300643272878
$ws.Range("B22").Value = 0.045490175617100448
$ws.Range("A23").Value = -0.040490175678356444
$ws.Range("B23").Value = 0.040097128268669735
$ws.Range("A24").Value = -0.020097128483864246
$ws.Range("B24").Value = 0.019999999781821209
$ws.Range("A25").Value = -0.027185735894931895
$ws.Range("B25").Value = 0.02711511693499169
$ws.Range("A26").Value = -0.024615116981193452
$ws.Range("B26").Value = 0.024525961184350464
$ws.Range("A27").Value = -0.022025961231693092
$ws.Range("B27").Value = 0.021509609644082328
$ws.Range("A28").Value = -0.019509609691507279
$ws.Range("B28").Value = 0.019173925009771331
$ws.Range("A29").Value = -0.012173925109563832
$ws.Range("B29").Value = 0.012086435171171139
$ws.Range("A30").Value = 0.047913564207433534
$ws.Range("B30").Value = -0.048225753042896446
$ws.Range("A31").Value = 0.055225752950903484
$ws.Range("B31").Value = -0.055327328462235315
$ws.Range("A32").Value = 0.065327328343645519
$ws.Range("B32").Value = -0.065483840255774695
